$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B33: should be a true numeric value (4) instead of text "4"
$ws.Range("B33").Value = 4

# Add new row 34 with the new annotation data
$ws.Range("A34").Value = "Sunsi Wu"
$ws.Range("B34").Value = "'3"
$ws.Range("B34").Style = "Normal"
$ws.Range("C34").Value = "无"
$ws.Range("D34").Value = "QSN"
$ws.Range("E34").Value = "RES"
$ws.Range("F34").Value = "2e6daeb6-f5b1-42e4-9927-e16202e5fb2e"
$ws.Range("G34").Value = "H1cWzoxA-_annotated.xlsx"
$ws.Range("H34").Value = 'For example ,when I use the cr dataset, "python sc_main.py --network_type exp_context_fusion --context_fusion_method wblock --model_dir_suffix training --dataset_type cr --gpu 0 " the result is not the 84.48 as the paper,I could only get 84.30 after several times.'
